$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 332935
$ws.Range("R2").Value = 6626957
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("A3").Value = 112164609
$ws.Range("Q3").Value = 332973
$ws.Range("R3").Value = 6627007
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
$ws.Range("A4").Value = 112164702
$ws.Range("B4").Value = 89369
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 5447
$ws.Range("F4").Value = "Vedticka"
$ws.Range("G4").Value = "Fuscoporia viticola"
$ws.Range("H4").Value = "(Schwein.) Murrill"
$ws.Range("L4").ClearContents()
$ws.Range("Q4").Value = 332980
$ws.Range("R4").Value = 6627033
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
$ws.Range("AC4").ClearContents()
$ws.Range("A5").Value = 112164565
$ws.Range("B5").Value = 92683
$ws.Range("E5").Value = 2362
$ws.Range("F5").Value = "Blek stjärnmossa"
$ws.Range("G5").Value = "Mnium stellare"
$ws.Range("H5").Value = "Hedw."
$ws.Range("L5").Value = ""
$ws.Range("Q5").Value = 332935
$ws.Range("R5").Value = 6626957
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()
$ws.Range("A6").Value = 112164673
$ws.Range("B6").Value = 93157
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 1078
$ws.Range("F6").Value = "Rundfjädermossa"
$ws.Range("G6").Value = "Neckera besseri"
$ws.Range("H6").Value = "(Lobarz.) Jur."
$ws.Range("Q6").Value = 332854
$ws.Range("R6").Value = 6626968
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()
$ws.Range("AC6").Value = "Under överhängande klippa"
$ws.Range("A7").Value = 112164607
$ws.Range("B7").Value = 93158
$ws.Range("E7").Value = 2667
$ws.Range("F7").Value = "Platt fjädermossa"
$ws.Range("G7").Value = "Neckera complanata"
$ws.Range("H7").Value = "(Hedw.) Huebener"
$ws.Range("Q7").Value = 332973
$ws.Range("R7").Value = 6627007
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()
$ws.Range("A8").Value = 112164661
$ws.Range("B8").Value = 89864
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 5467
$ws.Range("F8").Value = "Kådvaxskinn"
$ws.Range("G8").Value = "Phlebia serialis"
$ws.Range("H8").Value = "(Fr.:Fr.) Donk"
$ws.Range("L8").ClearContents()
$ws.Range("Q8").Value = 332865
$ws.Range("R8").Value = 6626972
$ws.Range("Z8").ClearContents()
$ws.Range("AB8").ClearContents()
$ws.Range("AC8").Value = "På granlåga"
$ws.Range("A9").Value = 112164579
$ws.Range("B9").Value = 93159
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 2666
$ws.Range("F9").Value = "Grov fjädermossa"
$ws.Range("G9").Value = "Neckera crispa"
$ws.Range("H9").Value = "Hedw."
$ws.Range("L9").Value = ""
$ws.Range("Q9").Value = 332923
$ws.Range("R9").Value = 6626955
$ws.Range("Z9").ClearContents()
$ws.Range("AB9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("Q10").Value = 332854
$ws.Range("R10").Value = 6626968
$ws.Range("Z10").ClearContents()
$ws.Range("AB10").ClearContents()
